# Auto-generated edit script: updates Leve profit-calculation tables
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# refreshed Universalis market-price data (columns H:N).

$wb = $excel.ActiveWorkbook

$sheetsData = @{
    "ALC" = @{
        19 = @{ writes = @{ "H" = 726.82355; "I" = 824.3333; "J" = 492.8; "K" = 824.3333; "L" = 492.8; "M" = -649.3333; "N" = -842.8 }; deletes = @() }
        28 = @{ writes = @{ "H" = 2442.182; "I" = 1620.625; "K" = 1620.625; "M" = -1135.625 }; deletes = @() }
        38 = @{ writes = @{ "H" = 255.2; "I" = 255.2; "K" = 765.5999999999999; "M" = -393.5999999999999 }; deletes = @() }
        40 = @{ writes = @{ "H" = 3368.5518; "I" = 1971.6111; "K" = 1971.6111; "M" = -1796.6111 }; deletes = @() }
        43 = @{ writes = @{ "H" = 3167.6667; "J" = 3002; "L" = 3002; "N" = -3140 }; deletes = @() }
        58 = @{ writes = @{ "H" = 631.25; "J" = 2000; "L" = 6000; "N" = -6300 }; deletes = @() }
        80 = @{ writes = @{ "H" = 1508.4445; "I" = 1307.8572; "J" = 1636.091; "K" = 3923.5716; "L" = 4908.272999999999; "M" = -2925.5716; "N" = -6904.272999999999 }; deletes = @() }
        83 = @{ writes = @{ "H" = 1508.4445; "I" = 1307.8572; "J" = 1636.091; "K" = 11770.7148; "L" = 14724.819; "M" = -6778.7148; "N" = -24708.819 }; deletes = @() }
        86 = @{ writes = @{ "H" = 1559.25; "I" = 1345.8334; "K" = 1345.8334; "M" = -222.8334 }; deletes = @() }
        89 = @{ writes = @{ "H" = 1559.25; "I" = 1345.8334; "K" = 6729.166999999999; "M" = -1113.166999999999 }; deletes = @() }
        98 = @{ writes = @{ "I" = 2417.4443; "J" = 3349; "K" = 2417.4443; "L" = 3349; "M" = -919.4443000000001; "N" = -6345 }; deletes = @() }
        100 = @{ writes = @{ "H" = 2471.3333; "I" = 2471.3333; "K" = 2471.3333; "M" = -1930.3333 }; deletes = @() }
        122 = @{ writes = @{ "I" = 2417.4443; "J" = 3349; "K" = 7252.3329; "L" = 10047; "M" = -4802.3329; "N" = -14947 }; deletes = @() }
        125 = @{ writes = @{ "H" = 595.5714; "I" = 567.6667; "J" = 603.1818; "K" = 5109.0003; "L" = 5428.6362; "M" = -2649.0003; "N" = -10348.6362 }; deletes = @() }
        135 = @{ writes = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }; deletes = @("M","N") }
    }
    "ARM" = @{
        32 = @{ writes = @{ "H" = 4314.6665; "I" = 2888.7273; "K" = 2888.7273; "M" = -2601.7273 }; deletes = @() }
        63 = @{ writes = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }; deletes = @("M","N") }
        66 = @{ writes = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }; deletes = @("M","N") }
        97 = @{ writes = @{ "H" = 3173.25; "I" = 897.6667; "J" = 10000; "K" = 897.6667; "L" = 10000; "M" = -401.6667; "N" = -10992 }; deletes = @() }
        102 = @{ writes = @{ "H" = 0; "I" = 0; "K" = 0 }; deletes = @("M") }
        122 = @{ writes = @{ "H" = 4754; "I" = 4754; "K" = 14262; "M" = -11812 }; deletes = @() }
    }
    "BSM" = @{
        20 = @{ writes = @{ "H" = 4125.75; "I" = 4125.75; "J" = 0; "K" = 4125.75; "L" = 0; "M" = -3878.75 }; deletes = @("N") }
        94 = @{ writes = @{ "H" = 2269.9285; "I" = 1815.0834; "K" = 1815.0834; "M" = -1364.0834 }; deletes = @() }
        99 = @{ writes = @{ "H" = 2049; "I" = 1123.75; "K" = 1123.75; "M" = 374.25 }; deletes = @() }
        105 = @{ writes = @{ "H" = 1074.25; "I" = 966; "K" = 966; "M" = 781 }; deletes = @() }
    }
    "CRP" = @{
        4 = @{ writes = @{ "H" = 50; "J" = 0; "L" = 0 }; deletes = @("N") }
        105 = @{ writes = @{ "H" = 1679.6666; "I" = 1199.3334; "K" = 1199.3334; "M" = 547.6666 }; deletes = @() }
        107 = @{ writes = @{ "H" = 660.35; "I" = 544.8889; "J" = 1699.5; "K" = 544.8889; "L" = 1699.5; "M" = 1375.1111; "N" = -5539.5 }; deletes = @() }
        132 = @{ writes = @{ "H" = 3151.8928; "I" = 2472.5454; "J" = 5642.8335; "K" = 7417.6362; "L" = 16928.5005; "M" = -4887.6362; "N" = -21988.5005 }; deletes = @() }
    }
    "CUL" = @{
        12 = @{ writes = @{ "H" = 55.18182; "J" = 61.555557; "L" = 184.666671; "N" = -530.666671 }; deletes = @() }
        49 = @{ writes = @{ "H" = 4249.5; "I" = 4499; "K" = 13497; "M" = -13341 }; deletes = @() }
        61 = @{ writes = @{ "H" = 188; "I" = 188; "J" = 0; "K" = 564; "L" = 0; "M" = -349 }; deletes = @("N") }
        98 = @{ writes = @{ "H" = 474.75; "I" = 399; "K" = 1197; "M" = 301 }; deletes = @() }
    }
    "GSM" = @{
        70 = @{ writes = @{ "H" = 9441; "I" = 6468.6665; "K" = 6468.6665; "M" = -6198.6665 }; deletes = @() }
        73 = @{ writes = @{ "H" = 9441; "I" = 6468.6665; "K" = 6468.6665; "M" = -5532.6665 }; deletes = @() }
        104 = @{ writes = @{ "H" = 99999.336; "J" = 99999.336; "L" = 99999.336; "N" = -106987.336 }; deletes = @() }
        132 = @{ writes = @{ "H" = 3266.818; "I" = 3172.2632; "J" = 3865.6667; "K" = 9516.7896; "L" = 11597.0001; "M" = -6986.7896; "N" = -16657.0001 }; deletes = @() }
    }
    "LTW" = @{
        39 = @{ writes = @{ "H" = 1750; "I" = 1000; "K" = 1000; "M" = -540 }; deletes = @() }
        55 = @{ writes = @{ "H" = 780.125; "I" = 340.2857; "J" = 1122.2222; "K" = 340.2857; "L" = 1122.2222; "M" = -167.2857; "N" = -1468.2222 }; deletes = @() }
        122 = @{ writes = @{ "H" = 2676.3635; "I" = 2041; "J" = 3438.8; "K" = 6123; "L" = 10316.4; "M" = -3673; "N" = -15216.4 }; deletes = @() }
        132 = @{ writes = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }; deletes = @("M","N") }
    }
    "WVR" = @{
        100 = @{ writes = @{ "H" = 278.57144; "I" = 278.57144; "K" = 557.14288; "M" = -16.14287999999999 }; deletes = @() }
        107 = @{ writes = @{ "H" = 444.36365; "I" = 476.8; "K" = 1430.4; "M" = 489.5999999999999 }; deletes = @() }
        113 = @{ writes = @{ "H" = 471.44446; "I" = 364.6; "J" = 605; "K" = 1093.8; "L" = 1815; "M" = 1076.2; "N" = -6155 }; deletes = @() }
    }
}

foreach ($sheetName in $sheetsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetsData[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $rowPlan = $rows[$rowNum]
        foreach ($col in $rowPlan.writes.Keys) {
            $ws.Range("$col$rowNum").Value = $rowPlan.writes[$col]
        }
        foreach ($col in $rowPlan.deletes) {
            $ws.Range("$col$rowNum").ClearContents()
        }
    }
}